$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Artfynd")

$ws.Range("B2").Value = 91829
$ws.Range("B3").Value = 79244
$ws.Range("B4").Value = 79244
$ws.Range("B5").Value = 91829
